$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats/styles from column F (the old column D, now shifted) into new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 37 and 79 are section-title rows with only a single label cell (B37/B79);
# the block PasteSpecial above spuriously created empty D/E cells there - remove them
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()

# Populate new column D and E with the new quarter data
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 3306800
$ws.Range("E8").Value2 = 3395100
$ws.Range("D9").Value2 = 2344700
$ws.Range("E9").Value2 = 2382500
$ws.Range("D10").Value2 = 962100
$ws.Range("E10").Value2 = 1012600
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 500
$ws.Range("E14").Value2 = -7000
$ws.Range("D15").Value2 = 171100
$ws.Range("E15").Value2 = 175300
$ws.Range("D17").Value2 = 2824600
$ws.Range("E17").Value2 = 2825400
$ws.Range("D18").Value2 = 482200
$ws.Range("E18").Value2 = 569700
$ws.Range("D20").Value2 = 14600
$ws.Range("E20").Value2 = 29000
$ws.Range("D21").Value2 = 667800
$ws.Range("E21").Value2 = 774000
$ws.Range("D22").Value2 = 3800
$ws.Range("E22").Value2 = 5400
$ws.Range("D23").Value2 = 493000
$ws.Range("E23").Value2 = 593300
$ws.Range("D24").Value2 = 135700
$ws.Range("E24").Value2 = 140700
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 357300
$ws.Range("E26").Value2 = 452600
$ws.Range("D27").Value2 = 337100
$ws.Range("E27").Value2 = 431700
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = "NA"
$ws.Range("E29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -14600
$ws.Range("E32").Value2 = -29000
$ws.Range("D33").Value2 = 337100
$ws.Range("E33").Value2 = 431700
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 337100
$ws.Range("E35").Value2 = 431700
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 887300
$ws.Range("E41").Value2 = 1229700
$ws.Range("D42").Value2 = 1181500
$ws.Range("E42").Value2 = 635400
$ws.Range("D43").Value2 = 639900
$ws.Range("E43").Value2 = 670800
$ws.Range("D44").Value2 = 290700
$ws.Range("E44").Value2 = 304900
$ws.Range("D45").Value2 = 289800
$ws.Range("E45").Value2 = 198000
$ws.Range("D46").Value2 = 3289200
$ws.Range("E46").Value2 = 3038900
$ws.Range("D47").Value2 = 119000
$ws.Range("E47").Value2 = 109500
$ws.Range("D48").Value2 = 1928200
$ws.Range("E48").Value2 = 1990000
$ws.Range("D49").Value2 = 24737000
$ws.Range("E49").Value2 = 24487100
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 513700
$ws.Range("E52").Value2 = 511800
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 30587000
$ws.Range("E54").Value2 = 30137400
$ws.Range("D57").Value2 = 233800
$ws.Range("E57").Value2 = 269500
$ws.Range("D58").Value2 = 1342000
$ws.Range("E58").Value2 = 1367100
$ws.Range("D59").Value2 = 3200400
$ws.Range("E59").Value2 = 3241000
$ws.Range("D60").Value2 = 4776100
$ws.Range("E60").Value2 = 4877600
$ws.Range("D61").Value2 = 13810800
$ws.Range("E61").Value2 = 13790100
$ws.Range("D62").Value2 = 2945800
$ws.Range("E62").Value2 = 2786500
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 21991300
$ws.Range("E66").Value2 = 21892600
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 5212800
$ws.Range("E72").Value2 = 4875700
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 8595700
$ws.Range("E76").Value2 = 8244800
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 337100
$ws.Range("E81").Value2 = 431700
$ws.Range("D83").Value2 = 171100
$ws.Range("E83").Value2 = 175300
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 502100
$ws.Range("E89").Value2 = 690400
$ws.Range("D91").Value2 = -111100
$ws.Range("E91").Value2 = -114700
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -817200
$ws.Range("E94").Value2 = -460600
$ws.Range("D96").Value2 = 0
$ws.Range("E96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -27000
$ws.Range("E100").Value2 = -7900
$ws.Range("D101").Value2 = 0
$ws.Range("E101").Value2 = 0
$ws.Range("D102").Value2 = -342200
$ws.Range("E102").Value2 = 221900

# Row 91 data correction: old G91/H91 values moved to I91/J91 after shift, but the shifted-in
# values were superseded by corrected figures in the source data
$ws.Range("I91").Value2 = -111700
$ws.Range("J91").Value2 = -97500
